$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J3").Value = 1.67
$ws.Range("K3").Value = 2.88
$ws.Range("N3").Value = 21
$ws.Range("U3").Value = 1.7
$ws.Range("V3").Value = 2.05
$ws.Range("Y3").Value = 9
$ws.Range("AH3").Value = 29
$ws.Range("AW3").Value = 9.5
$ws.Range("W4").Value = 5.5
$ws.Range("AH4").Value = 11
$ws.Range("AI4").Value = 23
$ws.Range("AO4").Value = 9
$ws.Range("AP4").Value = 23
$ws.Range("BC4").Value = 126
$ws.Range("M5").Value = 1.1
$ws.Range("N5").Value = 7
$ws.Range("O5").Value = 1.5
$ws.Range("P5").Value = 2.5
$ws.Range("Q5").Value = 2.6
$ws.Range("R5").Value = 1.48
$ws.Range("T6").Value = 2.08
$ws.Range("S10").Value = 1.44
$ws.Range("T10").Value = 2.63
$ws.Range("AF10").Value = 51
$ws.Range("AT10").Value = 2.63
$ws.Range("G13").Value = 2.55
$ws.Range("H13").Value = 3.6
$ws.Range("I13").Value = 2.45
$ws.Range("K13").Value = 2.3
$ws.Range("L13").Value = 2.92
$ws.Range("P13").Value = 4.6
$ws.Range("S13").Value = 1.29
$ws.Range("T13").Value = 3.3
$ws.Range("X13").Value = 16.5
$ws.Range("AD13").Value = 7.6
$ws.Range("AH13").Value = 13
$ws.Range("AP13").Value = 16.5
$ws.Range("AT13").Value = 3.3
$ws.Range("AX13").Value = 12.5
$ws.Range("AY13").Value = 16
$ws.Range("BA13").Value = 60
$ws.Range("H27").Value = 3.1
$ws.Range("I27").Value = 2.7
$ws.Range("L27").Value = 3.6
$ws.Range("M27").Value = 1.11
$ws.Range("N27").Value = 6.5
$ws.Range("X27").Value = 12
$ws.Range("Z27").Value = 29
$ws.Range("AL27").Value = 26
$ws.Range("AW27").Value = 4.5
$ws.Range("M32").Value = 1.03
$ws.Range("N32").Value = 17
$ws.Range("Q32").Value = 1.5
$ws.Range("R32").Value = 2.5
$ws.Range("N36").Value = 21
$ws.Range("O36").Value = 1.13
$ws.Range("P36").Value = 6
$ws.Range("Q36").Value = 1.44
$ws.Range("R36").Value = 2.7
$ws.Range("W38").Value = 7.5
$ws.Range("AC38").Value = 9.5
$ws.Range("AG38").Value = 251
$ws.Range("AH38").Value = 9.5
$ws.Range("AY38").Value = 26
$ws.Range("BC42").Value = 126
$ws.Range("G47").Value = 3.6
$ws.Range("I47").Value = 1.85
$ws.Range("K47").Value = 2.2
$ws.Range("L47").Value = 2.6
$ws.Range("M47").Value = 1.05
$ws.Range("N47").Value = 11
$ws.Range("Q47").Value = 2
$ws.Range("R47").Value = 1.85
$ws.Range("U47").Value = 1.8
$ws.Range("V47").Value = 1.91
$ws.Range("W47").Value = 11
$ws.Range("AA47").Value = 34
$ws.Range("AC47").Value = 10
$ws.Range("AD47").Value = 7
$ws.Range("AG47").Value = 251
$ws.Range("AH47").Value = 7.5
$ws.Range("AL47").Value = 15
$ws.Range("AM47").Value = 26
$ws.Range("AR47").Value = 81
$ws.Range("AX47").Value = 10
$ws.Range("AZ47").Value = 34
$ws.Range("J48").Value = 3.4
$ws.Range("AK48").Value = 23
$ws.Range("J50").Value = 2.25
$ws.Range("AA50").Value = 12
$ws.Range("AO50").Value = 8.5
$ws.Range("AU50").Value = 7.5
$ws.Range("G53").Value = 5.5
$ws.Range("H53").Value = 3.8
$ws.Range("I53").Value = 1.53
$ws.Range("AA53").Value = 41
$ws.Range("AD53").Value = 7.5
$ws.Range("AU53").Value = 8.5
$ws.Range("AV53").Value = 51
$ws.Range("AW53").Value = 3.5
$ws.Range("AZ53").Value = 26
$ws.Range("G59").Value = 2.05
$ws.Range("I59").Value = 3.25
$ws.Range("N59").Value = 21
$ws.Range("Z59").Value = 21
$ws.Range("AA59").Value = 15
$ws.Range("AC59").Value = 21
$ws.Range("AD59").Value = 8
$ws.Range("AE59").Value = 11
$ws.Range("AF59").Value = 29
$ws.Range("AH59").Value = 17
$ws.Range("AI59").Value = 21
$ws.Range("AJ59").Value = 12
$ws.Range("AK59").Value = 34
$ws.Range("AL59").Value = 21
$ws.Range("AW59").Value = 5.5
$ws.Range("AX59").Value = 15
$ws.Range("U67").Value = 1.54
$ws.Range("I68").Value = 3.9
$ws.Range("L68").Value = 4
$ws.Range("U68").Value = 1.54
$ws.Range("W68").Value = 10
$ws.Range("AO68").Value = 9.5
$ws.Range("AX68").Value = 19
$ws.Range("AZ68").Value = 51
$ws.Range("U69").Value = 1.54
$ws.Range("U70").Value = 1.47
$ws.Range("Q73").Value = 1.67
$ws.Range("R73").Value = 2.15
$ws.Range("I77").Value = 1.9
$ws.Range("J77").Value = 4.75
$ws.Range("L77").Value = 2.63
$ws.Range("N77").Value = 8
$ws.Range("AD77").Value = 6
$ws.Range("AF77").Value = 51
$ws.Range("AO77").Value = 23
$ws.Range("AT77").Value = 2.5
$ws.Range("AX77").Value = 11
$ws.Range("AZ77").Value = 41
$ws.Range("BA77").Value = 67
$ws.Range("G80").Value = 3.2
$ws.Range("I80").Value = 2.45
$ws.Range("AJ80").Value = 11
$ws.Range("AO80").Value = 19
$ws.Range("G81").Value = 4.33
$ws.Range("H81").Value = 3.5
$ws.Range("I81").Value = 1.8
$ws.Range("L81").Value = 2.4
$ws.Range("Z81").Value = 41
$ws.Range("AC81").Value = 11
$ws.Range("AD81").Value = 6.5
$ws.Range("AH81").Value = 7.5
$ws.Range("AL81").Value = 15
$ws.Range("AZ81").Value = 34
$ws.Range("Q86").Value = 2.25
$ws.Range("R86").Value = 1.62
$ws.Range("N89").Value = 17
$ws.Range("M97").Value = 1.03
$ws.Range("N97").Value = 17
$ws.Range("G101").Value = 13.5
$ws.Range("H101").Value = 5.1
$ws.Range("J101").Value = 10.75
$ws.Range("K101").Value = 2.35
$ws.Range("O101").Value = 1.23
$ws.Range("P101").Value = 3.35
$ws.Range("U101").Value = 2.42
$ws.Range("V101").Value = 1.44
$ws.Range("W101").Value = 28
$ws.Range("X101").Value = 120
$ws.Range("Z101").Value = 700
$ws.Range("AA101").Value = 300
$ws.Range("AC101").Value = 10.5
$ws.Range("AD101").Value = 11.25
$ws.Range("AE101").Value = 35
$ws.Range("AK101").Value = 6.5
$ws.Range("AL101").Value = 12
$ws.Range("AN101").Value = 12
$ws.Range("AT101").Value = 2.62
$ws.Range("AU101").Value = 10.75
$ws.Range("AX101").Value = 5.2
$ws.Range("BA101").Value = 50
$ws.Range("BB101").Value = 350
